# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the (Fecha, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Precio $/Kg) tuples across the
# data rows (2-19) of the sheet. Rows 8 and 18 are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @{
    2  = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    3  = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    4  = @{ D = 44922; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    5  = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    6  = @{ D = 44894; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    7  = @{ D = 44895; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    9  = @{ D = 44943; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    10 = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    11 = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    12 = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    13 = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
    14 = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    15 = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    16 = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    17 = @{ D = 44915; J = 50; K = 18000; L = 18000; M = 18000; P = 1385 }
    19 = @{ D = 44930; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
